$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - motor eléctrico
$ws.Range("A3").Value = "motor eléctrico"
$ws.Range("B3").Value = "simotics"
$ws.Range("C3").Value = "1hp a 440v"
$ws.Range("D3").Value = 3381
$ws.Range("E3").Value = 3381
$ws.Range("F3").Value = "simotics"

# Row 4 - cable 18AWG
$ws.Range("A4").Value = "cable 18AWG"
$ws.Range("B4").Value = "steren"
$ws.Range("C4").Value = "60 metros"
$ws.Range("D4").Value = 12
$ws.Range("E4").Value = 720
$ws.Range("F4").Value = "steren"

# Row 5 - pulsador de marcha
$ws.Range("A5").Value = "pulsador de marcha"

# Row 6 - pulsador de paro
$ws.Range("A6").Value = "pulsador de paro"

$ws.Range("A6").Select()
